$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-21, columns B (Matriz), D (Autovalor), E (Iteracoes), F (Tempo), G (Ordem)
# Values below reflect the re-synchronized run results (re-ordering + new timings)

$ws.Cells.Item(2, 2).Value = "can_187.mtx"
$ws.Cells.Item(2, 4).Value = 8.397620475292884
$ws.Cells.Item(2, 5).Value = 58
$ws.Cells.Item(2, 6).Value = 0.00758671760559082
$ws.Cells.Item(2, 7).Value = 187

$ws.Cells.Item(3, 2).Value = "can_187.mtx"
$ws.Cells.Item(3, 4).Value = 8.39929552537935
$ws.Cells.Item(3, 5).Value = 24
$ws.Cells.Item(3, 6).Value = 0.0009226799011230469
$ws.Cells.Item(3, 7).Value = 187

$ws.Cells.Item(4, 2).Value = "can_229.mtx"
$ws.Cells.Item(4, 4).Value = 8.696410230685041
$ws.Cells.Item(4, 5).Value = 13
$ws.Cells.Item(4, 6).Value = 0.0006277561187744141
$ws.Cells.Item(4, 7).Value = 229

$ws.Cells.Item(5, 2).Value = "can_229.mtx"
$ws.Cells.Item(5, 4).Value = 8.696434096273036
$ws.Cells.Item(5, 5).Value = 9
$ws.Cells.Item(5, 6).Value = 0.0004148483276367188
$ws.Cells.Item(5, 7).Value = 229

$ws.Cells.Item(6, 2).Value = "can_161.mtx"
$ws.Cells.Item(6, 4).Value = 8.821261565896524
$ws.Cells.Item(6, 5).Value = 21
$ws.Cells.Item(6, 6).Value = 0.0007724761962890625
$ws.Cells.Item(6, 7).Value = 161

$ws.Cells.Item(7, 2).Value = "can_161.mtx"
$ws.Cells.Item(7, 4).Value = 8.82129418362558
$ws.Cells.Item(7, 5).Value = 14
$ws.Cells.Item(7, 6).Value = 0.0006036758422851562
$ws.Cells.Item(7, 7).Value = 161

$ws.Cells.Item(8, 2).Value = "can_292.mtx"
$ws.Cells.Item(8, 4).Value = 12.16372151100421
$ws.Cells.Item(8, 5).Value = 11
$ws.Cells.Item(8, 6).Value = 0.0008475780487060547
$ws.Cells.Item(8, 7).Value = 292

$ws.Cells.Item(9, 2).Value = "can_292.mtx"
$ws.Cells.Item(9, 4).Value = 12.16373315524623
$ws.Cells.Item(9, 5).Value = 9
$ws.Cells.Item(9, 6).Value = 0.0005564689636230469
$ws.Cells.Item(9, 7).Value = 292

$ws.Cells.Item(10, 2).Value = "can_1054.mtx"
$ws.Cells.Item(10, 4).Value = 14.84373761133698
$ws.Cells.Item(10, 5).Value = 42
$ws.Cells.Item(10, 6).Value = 0.01749372482299805
$ws.Cells.Item(10, 7).Value = 1054

$ws.Cells.Item(11, 2).Value = "can_1054.mtx"
$ws.Cells.Item(11, 4).Value = 14.84457329193234
$ws.Cells.Item(11, 5).Value = 28
$ws.Cells.Item(11, 6).Value = 0.01064252853393555
$ws.Cells.Item(11, 7).Value = 1054

$ws.Cells.Item(12, 2).Value = "can_445.mtx"
$ws.Cells.Item(12, 4).Value = 8.950130371605427
$ws.Cells.Item(12, 5).Value = 43
$ws.Cells.Item(12, 6).Value = 0.002956867218017578
$ws.Cells.Item(12, 7).Value = 445

$ws.Cells.Item(13, 2).Value = "can_445.mtx"
$ws.Cells.Item(13, 4).Value = 8.950544673990166
$ws.Cells.Item(13, 5).Value = 20
$ws.Cells.Item(13, 6).Value = 0.0013580322265625
$ws.Cells.Item(13, 7).Value = 445

$ws.Cells.Item(14, 2).Value = "can_256.mtx"
$ws.Cells.Item(14, 4).Value = 16.03687867071401
$ws.Cells.Item(14, 5).Value = 16
$ws.Cells.Item(14, 6).Value = 0.0009183883666992188
$ws.Cells.Item(14, 7).Value = 256

$ws.Cells.Item(15, 2).Value = "can_256.mtx"
$ws.Cells.Item(15, 4).Value = 16.03715130336158
$ws.Cells.Item(15, 5).Value = 14
$ws.Cells.Item(15, 6).Value = 0.0008599758148193359
$ws.Cells.Item(15, 7).Value = 256

$ws.Cells.Item(16, 2).Value = "can_268.mtx"
$ws.Cells.Item(16, 4).Value = 14.45866874440652
$ws.Cells.Item(16, 5).Value = 18
$ws.Cells.Item(16, 6).Value = 0.001137256622314453
$ws.Cells.Item(16, 7).Value = 268

$ws.Cells.Item(17, 2).Value = "can_268.mtx"
$ws.Cells.Item(17, 4).Value = 14.45863188310433
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = 0.0005848407745361328
$ws.Cells.Item(17, 7).Value = 268

$ws.Cells.Item(18, 2).Value = "can_715.mtx"
$ws.Cells.Item(18, 4).Value = 15.30932234687733
$ws.Cells.Item(18, 5).Value = 48
$ws.Cells.Item(18, 6).Value = 0.006733417510986328
$ws.Cells.Item(18, 7).Value = 715

$ws.Cells.Item(19, 2).Value = "can_715.mtx"
$ws.Cells.Item(19, 4).Value = 15.31046671178859
$ws.Cells.Item(19, 5).Value = 29
$ws.Cells.Item(19, 6).Value = 0.004374027252197266
$ws.Cells.Item(19, 7).Value = 715

$ws.Cells.Item(20, 2).Value = "can_634.mtx"
$ws.Cells.Item(20, 4).Value = 13.85636022094015
$ws.Cells.Item(20, 5).Value = 74
$ws.Cells.Item(20, 6).Value = 0.008114337921142578
$ws.Cells.Item(20, 7).Value = 634

$ws.Cells.Item(21, 2).Value = "can_634.mtx"
$ws.Cells.Item(21, 4).Value = 13.86178904213738
$ws.Cells.Item(21, 5).Value = 56
$ws.Cells.Item(21, 6).Value = 0.006328582763671875
$ws.Cells.Item(21, 7).Value = 634
